# Update Betfair Back/Lay odds values in Sheet1 to reflect the latest
# market prices, per the daily data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 3.85
$ws.Range("H2").Value = 1.87
$ws.Range("I2").Value = 1.98
$ws.Range("J2").Value = 3.95

# Row 3
$ws.Range("I3").Value = 1.91
$ws.Range("J3").Value = 3.9
$ws.Range("Q3").Value = 1.81
$ws.Range("AG3").Value = 20
$ws.Range("AK3").Value = 210

# Row 4
$ws.Range("Q4").Value = 1.63
$ws.Range("U4").Value = 2.7
$ws.Range("AK4").Value = 32

# Row 5
$ws.Range("Q5").Value = 1.74

# Row 6
$ws.Range("P6").Value = 1.07

# Row 7
$ws.Range("F7").Value = 2.34
$ws.Range("G7").Value = 3.35
$ws.Range("H7").Value = 2.42
$ws.Range("I7").Value = 3.5
$ws.Range("J7").Value = 2.74
$ws.Range("K7").Value = 3.8

# Row 9
$ws.Range("J9").Value = 2.9
$ws.Range("K9").Value = 3.3
$ws.Range("P9").Value = 1.44
$ws.Range("Q9").Value = 2.88

# Row 10
$ws.Range("G10").Value = 3.1
$ws.Range("J10").Value = 2.72

# Row 11
$ws.Range("F11").Value = 1.92
$ws.Range("G11").Value = 2.14
$ws.Range("J11").Value = 2.8
$ws.Range("K11").Value = 3.65
$ws.Range("P11").Value = 1.08
$ws.Range("Q11").Value = 1.01

$wb.Save()
